$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 409.8889
$ws.Range("I2").Value = 409.8889
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 409.8889
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -296.8889
$ws.Range("N2").ClearContents()

# Row 17
$ws.Range("H17").Value = 281300.47
$ws.Range("I17").Value = 40000
$ws.Range("J17").Value = 288194.78
$ws.Range("K17").Value = 120000
$ws.Range("L17").Value = 864584.3400000001
$ws.Range("M17").Value = -119832
$ws.Range("N17").Value = -864920.3400000001

# Row 21
$ws.Range("H21").Value = 9665.888999999999
$ws.Range("I21").Value = 9999
$ws.Range("J21").Value = 9399.4
$ws.Range("K21").Value = 9999
$ws.Range("L21").Value = 9399.4
$ws.Range("M21").Value = -9531

# Row 23
$ws.Range("H23").Value = 9665.888999999999
$ws.Range("I23").Value = 9999
$ws.Range("J23").Value = 9399.4
$ws.Range("K23").Value = 9999
$ws.Range("L23").Value = 9399.4
$ws.Range("M23").Value = -9765

# Row 32
$ws.Range("H32").Value = 4716
$ws.Range("I32").Value = 3696
$ws.Range("J32").Value = 5282.6665
$ws.Range("K32").Value = 3696
$ws.Range("L32").Value = 5282.6665
$ws.Range("M32").Value = -3370

# Row 80
$ws.Range("H80").Value = 914917.1
$ws.Range("I80").Value = 3804694.8
$ws.Range("J80").Value = 2355.7368
$ws.Range("K80").Value = 11414084.4
$ws.Range("L80").Value = 7067.2104
$ws.Range("M80").Value = -11413086.4
$ws.Range("N80").Value = -9063.2104

# Row 83
$ws.Range("H83").Value = 914917.1
$ws.Range("I83").Value = 3804694.8
$ws.Range("J83").Value = 2355.7368
$ws.Range("K83").Value = 34242253.2
$ws.Range("L83").Value = 21201.6312
$ws.Range("M83").Value = -34237261.2
$ws.Range("N83").Value = -31185.6312

# Row 92
$ws.Range("H92").Value = 59457.234
$ws.Range("I92").Value = 756.6667
$ws.Range("J92").Value = 200338.6
$ws.Range("K92").Value = 756.6667
$ws.Range("L92").Value = 200338.6
$ws.Range("M92").Value = 491.3333

# Row 112
$ws.Range("H112").Value = 2747.7334
$ws.Range("I112").Value = 4767.6
$ws.Range("J112").Value = 1737.8
$ws.Range("K112").Value = 14302.8
$ws.Range("L112").Value = 5213.4
$ws.Range("M112").Value = -13194.8
$ws.Range("N112").Value = -7429.4

# Row 125
$ws.Range("H125").Value = 1379.2142
$ws.Range("I125").Value = 1197.75
$ws.Range("J125").Value = 2468
$ws.Range("K125").Value = 10779.75
$ws.Range("L125").Value = 22212
$ws.Range("M125").Value = -8319.75
$ws.Range("N125").Value = -27132

# Row 135
$ws.Range("H135").Value = 58825076
$ws.Range("I135").Value = 125000540
$ws.Range("J135").Value = 2440.111
$ws.Range("K135").Value = 1125004860
$ws.Range("L135").Value = 21960.999
$ws.Range("M135").Value = -1125002325

# Row 138
$ws.Range("H138").Value = 2372.394
$ws.Range("I138").Value = 1497
$ws.Range("J138").Value = 3422.8667
$ws.Range("K138").Value = 4491
$ws.Range("L138").Value = 10268.6001
$ws.Range("M138").Value = 649
$ws.Range("N138").Value = -20548.6001


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1489.1273
$ws.Range("I32").Value = 1479.6666
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 1479.6666
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -1192.6666

# Row 35
$ws.Range("H35").Value = 1179
$ws.Range("I35").Value = 1179
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1179
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -773
$ws.Range("N35").ClearContents()

# Row 45
$ws.Range("H45").Value = 696038.25
$ws.Range("I45").Value = 794900.9
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 794900.9
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -794523.9

# Row 74
$ws.Range("H74").Value = 37041452
$ws.Range("I74").Value = 43482804
$ws.Range("J74").Value = 3693.5
$ws.Range("K74").Value = 43482804
$ws.Range("L74").Value = 3693.5
$ws.Range("M74").Value = -43481930

# Row 77
$ws.Range("H77").Value = 37041452
$ws.Range("I77").Value = 43482804
$ws.Range("J77").Value = 3693.5
$ws.Range("K77").Value = 217414020
$ws.Range("L77").Value = 18467.5
$ws.Range("M77").Value = -217409652


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 2924.2222
$ws.Range("I15").Value = 3119
$ws.Range("J15").Value = 2534.6667
$ws.Range("K15").Value = 3119
$ws.Range("L15").Value = 2534.6667
$ws.Range("M15").Value = -2949
$ws.Range("N15").Value = -2874.6667

# Row 18
$ws.Range("H18").Value = 5409.222
$ws.Range("I18").Value = 2500
$ws.Range("J18").Value = 11227.667
$ws.Range("K18").Value = 2500
$ws.Range("L18").Value = 11227.667
$ws.Range("M18").Value = -2270
$ws.Range("N18").Value = -11687.667

# Row 22
$ws.Range("H22").Value = 953.5
$ws.Range("I22").Value = 805.25
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 805.25
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -455.25
$ws.Range("N22").Value = -1950

# Row 48
$ws.Range("H48").Value = 200046
$ws.Range("I48").Value = 200046
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 200046
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -199570
$ws.Range("N48").ClearContents()

# Row 58
$ws.Range("H58").Value = 2936.1667
$ws.Range("I58").Value = 1786.25
$ws.Range("J58").Value = 5236
$ws.Range("K58").Value = 1786.25
$ws.Range("L58").Value = 5236
$ws.Range("M58").Value = -1583.25
$ws.Range("N58").Value = -5642

# Row 99
$ws.Range("H99").Value = 2212.6924
$ws.Range("I99").Value = 1838.875
$ws.Range("J99").Value = 2810.8
$ws.Range("K99").Value = 1838.875
$ws.Range("L99").Value = 2810.8
$ws.Range("M99").Value = -340.875

# Row 122
$ws.Range("H122").Value = 2004.8667
$ws.Range("I122").Value = 2004.8667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6014.6001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3564.6001

# Row 126
$ws.Range("H126").Value = 2212.6924
$ws.Range("I126").Value = 1838.875
$ws.Range("J126").Value = 2810.8
$ws.Range("K126").Value = 5516.625
$ws.Range("L126").Value = 8432.400000000001
$ws.Range("M126").Value = -3046.625

# Row 136
$ws.Range("H136").Value = 2936.1667
$ws.Range("I136").Value = 1786.25
$ws.Range("J136").Value = 5236
$ws.Range("K136").Value = 5358.75
$ws.Range("L136").Value = 15708
$ws.Range("M136").Value = -2808.75
$ws.Range("N136").Value = -20808


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 9860403
$ws.Range("I4").Value = 15090671
$ws.Range("J4").Value = 3357.6155
$ws.Range("K4").Value = 45272013
$ws.Range("L4").Value = 10072.8465
$ws.Range("M4").Value = -45271901
$ws.Range("N4").Value = -10296.8465

# Row 34
$ws.Range("H34").Value = 434.42856
$ws.Range("I34").Value = 81.666664
$ws.Range("J34").Value = 699
$ws.Range("K34").Value = 244.999992
$ws.Range("L34").Value = 2097
$ws.Range("M34").Value = -160.999992
$ws.Range("N34").Value = -2265

# Row 109
$ws.Range("H109").Value = 1411.8572
$ws.Range("I109").Value = 813.8333
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 2441.4999
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -1401.4999

# Row 132
$ws.Range("H132").Value = 1378
$ws.Range("I132").Value = 1269.8
$ws.Range("J132").Value = 1594.4
$ws.Range("K132").Value = 11428.2
$ws.Range("L132").Value = 14349.6
$ws.Range("M132").Value = -8898.199999999999


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 126
$ws.Range("H126").Value = 7054.087
$ws.Range("I126").Value = 8201.875
$ws.Range("J126").Value = 4430.5713
$ws.Range("K126").Value = 24605.625
$ws.Range("L126").Value = 13291.7139
$ws.Range("M126").Value = -22135.625


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1564.2941
$ws.Range("I46").Value = 1269.5
$ws.Range("J46").Value = 1985.4286
$ws.Range("K46").Value = 1269.5
$ws.Range("L46").Value = 1985.4286
$ws.Range("M46").Value = -1081.5
$ws.Range("N46").Value = -2361.4286

# Row 50
$ws.Range("H50").Value = 25849.375
$ws.Range("I50").Value = 23500
$ws.Range("J50").Value = 28198.75
$ws.Range("K50").Value = 23500
$ws.Range("L50").Value = 28198.75
$ws.Range("M50").Value = -22863
$ws.Range("N50").Value = -29472.75

# Row 132
$ws.Range("H132").Value = 4153.1333
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 10683.923
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 32051.769
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -37111.769

# Row 136
$ws.Range("H136").Value = 1765.1666
$ws.Range("I136").Value = 1639.1351
$ws.Range("J136").Value = 2697.8
$ws.Range("K136").Value = 4917.4053
$ws.Range("L136").Value = 8093.400000000001
$ws.Range("M136").Value = -2367.4053


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 3046.25
$ws.Range("I100").Value = 4482.4
$ws.Range("J100").Value = 652.6667
$ws.Range("K100").Value = 8964.799999999999
$ws.Range("L100").Value = 1305.3334
$ws.Range("M100").Value = -8423.799999999999

# Row 103
$ws.Range("H103").Value = 25981.25
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 25981.25
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 25981.25
$ws.Range("N103").Value = -28325.25

# Row 119
$ws.Range("H119").Value = 43999.7
$ws.Range("I119").Value = 52500
$ws.Range("J119").Value = 43055.223
$ws.Range("K119").Value = 52500
$ws.Range("L119").Value = 43055.223
$ws.Range("M119").Value = -47662
$ws.Range("N119").Value = -52731.223

# Row 122
$ws.Range("H122").Value = 2217.92
$ws.Range("I122").Value = 2218.1052
$ws.Range("J122").Value = 2217.3333
$ws.Range("K122").Value = 6654.3156
$ws.Range("L122").Value = 6651.999899999999
$ws.Range("M122").Value = -4204.3156

# Row 126
$ws.Range("H126").Value = 1662.55
$ws.Range("I126").Value = 1393.4667
$ws.Range("J126").Value = 2469.8
$ws.Range("K126").Value = 4180.4001
$ws.Range("L126").Value = 7409.400000000001
$ws.Range("M126").Value = -1710.4001
$ws.Range("N126").Value = -12349.4

# Row 132
$ws.Range("H132").Value = 4474.846
$ws.Range("I132").Value = 3540.3333
$ws.Range("J132").Value = 8399.799999999999
$ws.Range("K132").Value = 10620.9999
$ws.Range("L132").Value = 25199.4
$ws.Range("M132").Value = -8090.999899999999

# Row 136
$ws.Range("H136").Value = 3573.5293
$ws.Range("I136").Value = 2483.4333
$ws.Range("J136").Value = 11749.25
$ws.Range("K136").Value = 7450.2999
$ws.Range("L136").Value = 35247.75
$ws.Range("M136").Value = -4900.2999

